$wb = $excel.ActiveWorkbook

# The workbook has duplicated data in the "展览" (Exhibitions) sheet and the
# "全部类型" (All types) sheet. Both sheets need the "想去人数" (want-to-go count)
# column F updated for the same four rows.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value = 10596
    $ws.Range("F7").Value = 105
    $ws.Range("F8").Value = 1307
    $ws.Range("F23").Value = 1700
}
